$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add description text to D8 (was empty inlineStr)
$ws.Cells.Item(8, 4).Value = "Cooling of Water on Irrigation of An Absorber"

# 2. Convert the "PARTS" column (E) from UPPERCASE to Title Case for rows 8-51
$partsMap = @{
    8  = "Shell Plate"
    9  = "2:1 Ellipsoidal Head"
    10 = "Lifting Lug"
    11 = "Lifting Lug"
    12 = "Saddle Plate"
    13 = "Saddle Baseplate"
    14 = "Saddle Baseplate"
    15 = "Saddle Baseplate"
    16 = "Flange"
    17 = "Flange"
    18 = "Neck"
    19 = "Neck"
    20 = "Flange"
    21 = "Flange"
    22 = "Neck"
    23 = "Flange"
    24 = "Neck"
    25 = "Neck"
    26 = "Flange"
    27 = "Neck"
    28 = "Neck"
    29 = "Blind Flange"
    30 = "Spiral Wound Gasket"
    31 = "Flange"
    32 = "Neck"
    33 = "Stud Bolt"
    34 = "Nuts & Washer"
    35 = "Bracket 1"
    36 = "Angle Bar"
    37 = "Bracket 2"
    38 = "Bracket 3"
    39 = "Grating"
    40 = "Angle Bar"
    41 = "Angle Bar"
    42 = "Flat Plate"
    43 = "Flat Plate"
    44 = "Angle Bar"
    45 = "Angle Bar"
    46 = "Side Rail"
    47 = "Ladder Bracket"
    48 = "Round Bar"
    49 = "Saddle Support Plate"
    50 = "Pad Plate"
    51 = "Pad Plate"
}
foreach ($row in $partsMap.Keys) {
    $ws.Cells.Item($row, 5).Value = $partsMap[$row]
}

# 3. Normalize "4 Bar.G" -> "4 Bar G" (column M) and "1 Bar.G" -> "1 Bar G" (column O) for rows 8-51
for ($row = 8; $row -le 51; $row++) {
    $ws.Cells.Item($row, 13).Value = "4 Bar G"
    $ws.Cells.Item($row, 15).Value = "1 Bar G"
}

# 4. Row 30 (Spiral Wound Gasket) material info corrections
$ws.Cells.Item(30, 8).Value = "Stainless Steel"
$ws.Cells.Item(30, 9).Value = "ASTM A182"
$ws.Cells.Item(30, 10).Value = "F304L"

# 5. Row 33 (Stud Bolt) grade correction
$ws.Cells.Item(33, 10).Value = "GR B8M"

# 6. Row 34 (Nuts & Washer) grade correction
$ws.Cells.Item(34, 10).Value = "GR 2H"

# 7. Row 39 (Grating) add spec
$ws.Cells.Item(39, 9).Value = "ASTM A36"

# 8. Remove rows 52-56 (Ladder, Handrail, Platform, Shell, Shell) -
#    this also shrinks the dimension and the A/B/C/D merged ranges from ..56 to ..51
$ws.Range("A52:O56").EntireRow.Delete()
